# Updated cryptos list - refresh Price (D) and Volume(1h) (E) columns
# for rows 2-51 on the active worksheet, matching the latest scrape.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @(
    @{ Row = 2;  D = "42.878.31";  E = "  -2.16%  " },
    @{ Row = 3;  D = "2.578.29";   E = "  +0.24%  " },
    @{ Row = 4;  D = $null;        E = "  -0.13%  " },
    @{ Row = 5;  D = "302.49";     E = "  +0.18%  " },
    @{ Row = 6;  D = "97.13";      E = "  +0.67%  " },
    @{ Row = 7;  D = "0.575";      E = "  -1.10%  " },
    @{ Row = 8;  D = $null;        E = "  -0.07%  " },
    @{ Row = 9;  D = "0.551";      E = "  -2.44%  " },
    @{ Row = 10; D = "36.63";      E = "  -1.49%  " },
    @{ Row = 11; D = $null;        E = "  -0.57%  " },
    @{ Row = 12; D = "7.74";       E = "  -1.15%  " },
    @{ Row = 13; D = $null;        E = "  +5.93%  " },
    @{ Row = 14; D = "2.586.79";   E = "  +0.32%  " },
    @{ Row = 15; D = $null;        E = "  -0.48%  " },
    @{ Row = 16; D = $null;        E = "  +0.11%  " },
    @{ Row = 17; D = "42.905.33";  E = "  -2.26%  " },
    @{ Row = 18; D = $null;        E = "  +3.12%  " },
    @{ Row = 19; D = $null;        E = "  +1.10%  " },
    @{ Row = 20; D = $null;        E = "  -0.91%  " },
    @{ Row = 21; D = "72.07";      E = "  -1.98%  " },
    @{ Row = 22; D = "254.86";     E = "  -4.27%  " },
    @{ Row = 23; D = $null;        E = "  +0.31%  " },
    @{ Row = 24; D = $null;        E = "  -5.17%  " },
    @{ Row = 25; D = "28.87";      E = "  -1.95%  " },
    @{ Row = 26; D = "1.00";       E = "  -0.15%  " },
    @{ Row = 27; D = $null;        E = "  +0.34%  " },
    @{ Row = 28; D = "37.86";      E = "  -1.36%  " },
    @{ Row = 29; D = $null;        E = "  -5.74%  " },
    @{ Row = 30; D = $null;        E = "  -2.69%  " },
    @{ Row = 31; D = "155.07";     E = "  +1.64%  " },
    @{ Row = 32; D = $null;        E = "  -1.49%  " },
    @{ Row = 33; D = "3.41";       E = "  -5.06%  " },
    @{ Row = 34; D = $null;        E = "  -2.08%  " },
    @{ Row = 35; D = "0.0805";     E = "  -1.31%  " },
    @{ Row = 36; D = "18.14";      E = "  +8.11%  " },
    @{ Row = 37; D = $null;        E = "  -3.24%  " },
    @{ Row = 38; D = $null;        E = "  -0.66%  " },
    @{ Row = 39; D = "23.13";      E = "  -3.46%  " },
    @{ Row = 40; D = $null;        E = "  -4.28%  " },
    @{ Row = 41; D = $null;        E = "  -1.98%  " },
    @{ Row = 42; D = $null;        E = "  -0.08%  " },
    @{ Row = 43; D = "2.06";       E = "  +26.07%  " },
    @{ Row = 44; D = "2.073.73";   E = "  +1.79%  " },
    @{ Row = 45; D = $null;        E = "  -0.09%  " },
    @{ Row = 46; D = $null;        E = "  +0.69%  " },
    @{ Row = 47; D = "85.51";      E = "  -2.90%  " },
    @{ Row = 48; D = "76.92";      E = "  +10.54%  " },
    @{ Row = 49; D = "106.80";     E = "  +1.19%  " },
    @{ Row = 50; D = "2.826.09";   E = "  -0.23%  " },
    @{ Row = 51; D = $null;        E = "  -0.11%  " }
)

foreach ($u in $updates) {
    $row = $u.Row

    if ($null -ne $u.D) {
        # Force the cell to stay a text value (the sheet stores prices like
        # "302.49" or "42.878.31" as plain inline strings, not numbers).
        # Without this, Excel's COM layer auto-coerces number-looking text
        # into a real numeric value when .Value is assigned.
        $dCell = $ws.Cells.Item($row, 4)
        $dCell.NumberFormat = "@"
        $dCell.Value = $u.D
        $dCell.Style = "Normal"
    }

    if ($null -ne $u.E) {
        $ws.Cells.Item($row, 5).Value = $u.E
    }
}
